$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.927.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''1.641.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = '''1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''216.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''0.5041'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = '''1.004'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''0.2567'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''0.06402'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = '''19.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = '''0.07740'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '''4.268'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''1.632.39'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.20%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''1.862.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''0.5461'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''0.0₅7920'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''64.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''25.912.38'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = '''1.004'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = '''203.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.19%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''4.385'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''9.925'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.48%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = '''5.982'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.06%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''1.005'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''1.939'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.92%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = '''141.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.73%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''0.1139'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.14%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = '''15.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = '''6.771'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.97%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = '''1.247'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = '''0.04953'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.83%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = '''3.284'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.26%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = '''3.192'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '''1.549'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = '''2.376'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.08%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = '''2.632'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.84%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = '''0.8941'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.31%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''1.162.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''0.5617'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.79%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''0.01568'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.27%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''1.005'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = '''5.652'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''0.8091'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.81%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = 'Quant'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = '''99.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = '''1.775.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = '''0.0₈118'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.09%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = '''0.4545'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = 'Aave'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = '''55.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = 'Frax'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = '''1.001'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = '''0.05059'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = '''1.002'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.19%  '
$ws.Range("E51").Style = "Normal"
